$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Values ----

# New header cell for column D ("Mês")
$ws.Range("D1").Value = "Mês"

# Row 2 (existing row): A2 becomes numeric 0, B2/C2 updated, D2 new
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Default"
$ws.Range("C2").Value = "Cypress"
$ws.Range("D2").Value = "Modelo 1"

# New rows 3-8
$data = @(
    @(1, "Mês 1", "Cypress", "Modelo 1"),
    @(2, "Mês 2", "Cypress", "Modelo 2"),
    @(3, "Mês 3", "Cypress", "Modelo 2"),
    @(4, "Mês 4", "Cypress", "Modelo 3"),
    @(5, "Mês 5", "Cypress", "Modelo 3"),
    @(6, "Mês 6", "Cypress", "Modelo 3")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---- Formatting ----

# Header row (A1:D1) - centered, thin border, blue fill, themed font color
$header = $ws.Range("A1:D1")
$header.HorizontalAlignment = -4108 # xlCenter
$header.Borders.LineStyle = 1
$header.Interior.Color = 16308937
$header.Interior.Pattern = 1
$header.Font.Name = "Arial"
$header.Font.Size = 10
$header.Font.ThemeColor = 0

# Data rows (A2:D8) - centered, thin border, no fill, themed font color
$bodyData = $ws.Range("A2:D8")
$bodyData.HorizontalAlignment = -4108 # xlCenter
$bodyData.Borders.LineStyle = 1
$bodyData.Font.Name = "Arial"
$bodyData.Font.Size = 10
$bodyData.Font.ThemeColor = 0

# Column B width (~16.38 chars; engine quantizes to nearest 1/6 step)
$ws.Columns.Item(2).ColumnWidth = 15.5
